# Commit: Vincular Desenvolvedor De Requisito a Release
# - Adds two new sheets: Requisito_Sprint and Desenvolvedor_Requisito_Sprint
# - Adds a new Projeto row, two new Requisito rows, rewrites Sprint data rows
# - Makes the last sheet (Desenvolvedor_Requisito_Sprint) the active tab

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force text storage so numeric-looking / date-looking / boolean-looking
    # strings ("30000", "01/01/2001", "false", "5", ...) are not silently
    # auto-converted by Excel's type inference.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
}

function Set-NumValue($cell, $val) {
    $cell.Value = $val
}

# ---------------------------------------------------------------------
# 1. Projeto sheet: add a new project row (row 7)
# ---------------------------------------------------------------------
$wsProjeto = $wb.Worksheets.Item("Projeto")
Set-NumValue  $wsProjeto.Cells.Item(7,1) 6
Set-TextValue $wsProjeto.Cells.Item(7,2) "projeto d"
Set-NumValue  $wsProjeto.Cells.Item(7,3) 2
Set-NumValue  $wsProjeto.Cells.Item(7,4) 30000
Set-TextValue $wsProjeto.Cells.Item(7,5) "30000"

# ---------------------------------------------------------------------
# 2. Requisito sheet: add two new requirement rows (rows 5 and 6)
# ---------------------------------------------------------------------
$wsRequisito = $wb.Worksheets.Item("Requisito")
Set-NumValue  $wsRequisito.Cells.Item(5,1) 4
Set-TextValue $wsRequisito.Cells.Item(5,2) "aaa"
Set-NumValue  $wsRequisito.Cells.Item(5,3) 3
Set-NumValue  $wsRequisito.Cells.Item(5,4) 2
Set-NumValue  $wsRequisito.Cells.Item(5,5) 1

Set-NumValue  $wsRequisito.Cells.Item(6,1) 5
Set-TextValue $wsRequisito.Cells.Item(6,2) "5"
Set-NumValue  $wsRequisito.Cells.Item(6,3) 2
Set-NumValue  $wsRequisito.Cells.Item(6,4) 2
Set-NumValue  $wsRequisito.Cells.Item(6,5) 1

# ---------------------------------------------------------------------
# 3. Sprint sheet: replace the 2 existing data rows with 5 new ones
# ---------------------------------------------------------------------
$wsSprint = $wb.Worksheets.Item("Sprint")

Set-NumValue  $wsSprint.Cells.Item(2,1) 1
Set-TextValue $wsSprint.Cells.Item(2,2) "AAABB"
Set-TextValue $wsSprint.Cells.Item(2,3) "11/03/2019"
Set-TextValue $wsSprint.Cells.Item(2,4) "15/04/2019"
Set-TextValue $wsSprint.Cells.Item(2,5) "Finalizada"

Set-NumValue  $wsSprint.Cells.Item(3,1) 2
Set-TextValue $wsSprint.Cells.Item(3,2) "VAAAI"
Set-TextValue $wsSprint.Cells.Item(3,3) "01/01/2001"
Set-TextValue $wsSprint.Cells.Item(3,4) "01/02/2002"
Set-TextValue $wsSprint.Cells.Item(3,5) "Finalizada"

Set-NumValue  $wsSprint.Cells.Item(4,1) 3
Set-TextValue $wsSprint.Cells.Item(4,2) "aaa"
Set-TextValue $wsSprint.Cells.Item(4,3) "2017"
Set-TextValue $wsSprint.Cells.Item(4,4) "2018"
Set-TextValue $wsSprint.Cells.Item(4,5) "Finalizada"

Set-NumValue  $wsSprint.Cells.Item(5,1) 4
Set-TextValue $wsSprint.Cells.Item(5,2) "a"
Set-TextValue $wsSprint.Cells.Item(5,3) "1"
Set-TextValue $wsSprint.Cells.Item(5,4) "2"
Set-TextValue $wsSprint.Cells.Item(5,5) "Finalizada"

Set-NumValue  $wsSprint.Cells.Item(6,1) 5
Set-TextValue $wsSprint.Cells.Item(6,2) "vv"
Set-TextValue $wsSprint.Cells.Item(6,3) "2018"
Set-TextValue $wsSprint.Cells.Item(6,4) "2019"
Set-TextValue $wsSprint.Cells.Item(6,5) "Em Andamento"

# ---------------------------------------------------------------------
# 4. Add new sheet "Requisito_Sprint" at the end
# ---------------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIdx)
$wsReqSprint = $wb.Worksheets.Add($null, $lastSheet)
$wsReqSprint.Name = "Requisito_Sprint"

Set-TextValue $wsReqSprint.Cells.Item(1,1) "ID"
Set-TextValue $wsReqSprint.Cells.Item(1,2) "Id Requisito"
Set-TextValue $wsReqSprint.Cells.Item(1,3) "Id Sprint"
Set-TextValue $wsReqSprint.Cells.Item(1,4) "Vinculou Desenvolvedor"
Set-TextValue $wsReqSprint.Cells.Item(1,5) "NívelImpactoAlterações"

$reqSprintRows = @(
    @(1, 1, 1, "false", "Baixo"),
    @(2, 3, 2, "false", "Baixo"),
    @(3, 1, 1, "false", "Alto"),
    @(4, 3, 1, "false", "Alto"),
    @(5, 1, 1, "false", "Médio"),
    @(6, 4, 5, "false", "Médio"),
    @(7, 2, 4, "false", "Médio"),
    @(8, 2, 5, "false", "Baixo"),
    @(9, 2, 4, "false", "Médio"),
    @(10, 3, 5, "false", "Alto"),
    @(11, 3, 5, "false", "Alto"),
    @(12, 2, 5, "false", "Baixo"),
    @(13, 1, 5, "false", "Médio"),
    @(14, 3, 5, "false", "Médio"),
    @(15, 5, 5, "false", "Médio")
)

$r = 2
foreach ($row in $reqSprintRows) {
    Set-NumValue  $wsReqSprint.Cells.Item($r,1) $row[0]
    Set-NumValue  $wsReqSprint.Cells.Item($r,2) $row[1]
    Set-NumValue  $wsReqSprint.Cells.Item($r,3) $row[2]
    Set-TextValue $wsReqSprint.Cells.Item($r,4) $row[3]
    Set-TextValue $wsReqSprint.Cells.Item($r,5) $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 5. Add new sheet "Desenvolvedor_Requisito_Sprint" at the end
# ---------------------------------------------------------------------
$lastIdx2 = $wb.Worksheets.Count
$lastSheet2 = $wb.Worksheets.Item($lastIdx2)
$wsDevReqSprint = $wb.Worksheets.Add($null, $lastSheet2)
$wsDevReqSprint.Name = "Desenvolvedor_Requisito_Sprint"

Set-TextValue $wsDevReqSprint.Cells.Item(1,1) "ID Requisito_Sprint"
Set-TextValue $wsDevReqSprint.Cells.Item(1,2) "Id Desenvolvedor"
Set-TextValue $wsDevReqSprint.Cells.Item(1,3) "Porcentagem"

$devReqSprintRows = @(
    @(1, 1, 60),
    @(1, 2, 30),
    @(1, 1, 10),
    @(3, 1, 50),
    @(3, 2, 50),
    @(6, 1, 60),
    @(6, 2, 40),
    @(2, 1, 50),
    @(8, 1, 50),
    @(10, 1, 50),
    @(10, 2, 25)
)

$r2 = 2
foreach ($row in $devReqSprintRows) {
    Set-NumValue $wsDevReqSprint.Cells.Item($r2,1) $row[0]
    Set-NumValue $wsDevReqSprint.Cells.Item($r2,2) $row[1]
    Set-NumValue $wsDevReqSprint.Cells.Item($r2,3) $row[2]
    $r2 = $r2 + 1
}

# ---------------------------------------------------------------------
# 6. Make Desenvolvedor_Requisito_Sprint (the last sheet) the active tab
# ---------------------------------------------------------------------
$wsDevReqSprint.Activate()
